# Weekly refresh of the "Hortaliza, Terminal La Palmera de La Serena - Albahaca"
# sheet: the Fecha (D), Volumen (J), Precio mínimo (K), Precio máximo (L),
# Precio promedio ponderado (M) and Precio $/Kg (P) columns are re-sampled for
# the existing rows (2-28), and 3 new observation rows (29-31) are appended.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: RowNum, Fecha(D), Volumen(J), PrecioMin(K), PrecioMax(L), PrecioProm(M), PrecioKg(P)
$updates = @(
    @(2, 44377, 600, 4000, 4500, 4250, 4250),
    @(3, 44356, 600, 3000, 3500, 3250, 3250),
    @(4, 44349, 560, 3000, 3500, 3250, 3250),
    @(5, 44445, 600, 4500, 5000, 4750, 4750),
    @(6, 44425, 900, 4500, 5000, 4750, 4750),
    @(7, 44162, 2000, 2800, 3000, 2900, 2900),
    @(8, 44165, 1000, 3000, 3500, 3250, 3250),
    @(9, 44169, 2400, 3000, 3500, 3250, 3250),
    @(10, 44418, 800, 4500, 5000, 4750, 4750),
    @(11, 44172, 760, 3000, 3500, 3250, 3250),
    @(12, 44407, 720, 4000, 4500, 4250, 4250),
    @(13, 44334, 760, 3000, 3500, 3250, 3250),
    @(14, 44365, 800, 3500, 4000, 3750, 3750),
    @(15, 44397, 800, 4000, 4500, 4250, 4250),
    @(16, 44335, 600, 3000, 3500, 3250, 3250),
    @(17, 44434, 600, 4500, 5000, 4750, 4750),
    @(18, 44420, 900, 4500, 5000, 4750, 4750),
    @(19, 44348, 700, 3000, 3500, 3250, 3250),
    @(20, 44427, 600, 4500, 5000, 4750, 4750),
    @(21, 44341, 700, 3000, 3500, 3250, 3250),
    @(22, 44176, 2000, 3000, 3500, 3250, 3250),
    @(23, 44441, 600, 4500, 5000, 4750, 4750),
    @(24, 44432, 900, 4500, 5000, 4750, 4750),
    @(25, 44342, 560, 3000, 3500, 3250, 3250),
    @(26, 44379, 800, 4000, 4500, 4250, 4250),
    @(27, 44315, 700, 2500, 3000, 2750, 2750),
    @(28, 44446, 800, 4500, 5000, 4750, 4750)
)

foreach ($row in $updates) {
    $r = $row[0]
    $ws.Cells.Item($r, 4).Value = $row[1]
    $ws.Cells.Item($r, 10).Value = $row[2]
    $ws.Cells.Item($r, 11).Value = $row[3]
    $ws.Cells.Item($r, 12).Value = $row[4]
    $ws.Cells.Item($r, 13).Value = $row[5]
    $ws.Cells.Item($r, 16).Value = $row[6]
}

# New rows appended at the bottom (29-31), same shape as the existing ones.
$newRows = @(
    @(29, 44411, 880, 4000, 4500, 4250, 4250),
    @(30, 44435, 1500, 4500, 5000, 4750, 4750),
    @(31, 44314, 800, 2500, 3000, 2750, 2750)
)

foreach ($row in $newRows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = 8
    $ws.Cells.Item($r, 2).Value = "Terminal La Palmera de La Serena"
    $ws.Cells.Item($r, 3).Value = "Coquimbo"
    $ws.Cells.Item($r, 4).Value = $row[1]
    $ws.Cells.Item($r, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($r, 5).Value = 4
    $ws.Cells.Item($r, 6).Value = 100112052
    $ws.Cells.Item($r, 7).Value = "Albahaca"
    $ws.Cells.Item($r, 8).Value = "Sin especificar"
    $ws.Cells.Item($r, 9).Value = "Primera"
    $ws.Cells.Item($r, 10).Value = $row[2]
    $ws.Cells.Item($r, 11).Value = $row[3]
    $ws.Cells.Item($r, 12).Value = $row[4]
    $ws.Cells.Item($r, 13).Value = $row[5]
    $ws.Cells.Item($r, 14).Value = "`$/paquete"
    $ws.Cells.Item($r, 15).Value = "Región de Arica y Parinacota"
    $ws.Cells.Item($r, 16).Value = $row[6]
    $ws.Cells.Item($r, 17).Value = 1
    $ws.Cells.Item($r, 18).Value = "Hortaliza"
}

Write-Output "Updated rows 2-28 and appended rows 29-31"
